$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.115
$ws.Range("E3").Value = 0.126
$ws.Range("E4").Value = 0.172
$ws.Range("G4").Value = 0.001
$ws.Range("D5").Value = 0.002
$ws.Range("E5").Value = 0.294
$ws.Range("F5").Value = 0.003
$ws.Range("G5").Value = 0.004
$ws.Range("H5").Value = 0.003
$ws.Range("D6").Value = 0.058
$ws.Range("E6").Value = 0.334
$ws.Range("F6").Value = 0.079
$ws.Range("G6").Value = 0.078
$ws.Range("H6").Value = 0.058
$ws.Range("I6").Value = 0.069
$ws.Range("D7").Value = 0.525
$ws.Range("E7").Value = 0.516
$ws.Range("F7").Value = 0.507
$ws.Range("G7").Value = 0.517
$ws.Range("H7").Value = 0.503
$ws.Range("I7").Value = 0.514
$ws.Range("D8").Value = 0.965
$ws.Range("E8").Value = 0.666
$ws.Range("F8").Value = 0.958
$ws.Range("H8").Value = 0.965
$ws.Range("I8").Value = 0.951
$ws.Range("E9").Value = 0.833
$ws.Range("E10").Value = 0.945
$ws.Range("E11").Value = 0.997
$ws.Range("F13").Value = 0.991
$ws.Range("H13").Value = 0.457
$ws.Range("D14").Value = 1
$ws.Range("F14").Value = 0.977
$ws.Range("H14").Value = 0.467
$ws.Range("D15").Value = 0.998
$ws.Range("F15").Value = 0.947
$ws.Range("G15").Value = 0.999
$ws.Range("H15").Value = 0.459
$ws.Range("I15").Value = 1
$ws.Range("D16").Value = 0.973
$ws.Range("E16").Value = 0.999
$ws.Range("F16").Value = 0.865
$ws.Range("G16").Value = 0.981
$ws.Range("H16").Value = 0.486
$ws.Range("I16").Value = 0.986
$ws.Range("D17").Value = 0.859
$ws.Range("E17").Value = 0.924
$ws.Range("F17").Value = 0.737
$ws.Range("H17").Value = 0.48
$ws.Range("I17").Value = 0.893
$ws.Range("E18").Value = 0.504
$ws.Range("F18").Value = 0.506
$ws.Range("G18").Value = 0.502
$ws.Range("H18").Value = 0.494
$ws.Range("I18").Value = 0.495
$ws.Range("D19").Value = 0.084
$ws.Range("E19").Value = 0.06
$ws.Range("F19").Value = 0.237
$ws.Range("G19").Value = 0.079
$ws.Range("H19").Value = 0.503
$ws.Range("I19").Value = 0.073
$ws.Range("D20").Value = 0.001
$ws.Range("E20").Value = 0.001
$ws.Range("F20").Value = 0.048
$ws.Range("G20").Value = 0.001
$ws.Range("H20").Value = 0.514
$ws.Range("I20").Value = 0
$ws.Range("F21").Value = 0.002
$ws.Range("H21").Value = 0.57
$ws.Range("F22").Value = 0.001
$ws.Range("H22").Value = 0.608
$ws.Range("H23").Value = 0.645
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0.032
$ws.Range("F24").Value = 0.093
$ws.Range("G24").Value = 0.003
$ws.Range("H24").Value = 0.627
$ws.Range("I24").Value = 0
$ws.Range("D25").Value = 0.003
$ws.Range("E25").Value = 0.065
$ws.Range("F25").Value = 0.135
$ws.Range("G25").Value = 0.005
$ws.Range("H25").Value = 0.615
$ws.Range("I25").Value = 0.005
$ws.Range("E26").Value = 0.094
$ws.Range("F26").Value = 0.184
$ws.Range("G26").Value = 0.023
$ws.Range("H26").Value = 0.597
$ws.Range("I26").Value = 0.016
$ws.Range("D27").Value = 0.049
$ws.Range("E27").Value = 0.164
$ws.Range("F27").Value = 0.23
$ws.Range("G27").Value = 0.073
$ws.Range("H27").Value = 0.548
$ws.Range("I27").Value = 0.06
$ws.Range("D28").Value = 0.185
$ws.Range("E28").Value = 0.289
$ws.Range("F28").Value = 0.335
$ws.Range("G28").Value = 0.196
$ws.Range("H28").Value = 0.534
$ws.Range("I28").Value = 0.185
$ws.Range("D29").Value = 0.481
$ws.Range("E29").Value = 0.494
$ws.Range("F29").Value = 0.5
$ws.Range("G29").Value = 0.487
$ws.Range("H29").Value = 0.474
$ws.Range("I29").Value = 0.478
$ws.Range("D30").Value = 0.879
$ws.Range("E30").Value = 0.746
$ws.Range("F30").Value = 0.693
$ws.Range("G30").Value = 0.866
$ws.Range("H30").Value = 0.486
$ws.Range("I30").Value = 0.864
$ws.Range("F31").Value = 0.903
$ws.Range("G31").Value = 0.995
$ws.Range("H31").Value = 0.35
$ws.Range("I31").Value = 0.996
$ws.Range("E32").Value = 0.997
$ws.Range("F32").Value = 0.988
$ws.Range("H32").Value = 0.27
$ws.Range("H33").Value = 0.145
$ws.Range("H34").Value = 0.034
$ws.Range("E35").Value = 0.051
$ws.Range("E36").Value = 0.095
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0.152
$ws.Range("I37").Value = 0.994
$ws.Range("E38").Value = 0.213
$ws.Range("F38").Value = 0.98
$ws.Range("G38").Value = 0.968
$ws.Range("H38").Value = 0.998
$ws.Range("I38").Value = 0.975
$ws.Range("E39").Value = 0.339
$ws.Range("F39").Value = 0.866
$ws.Range("H39").Value = 0.94
$ws.Range("I39").Value = 0.844
$ws.Range("D40").Value = 0.496
$ws.Range("E40").Value = 0.478
$ws.Range("F40").Value = 0.499
$ws.Range("G40").Value = 0.476
$ws.Range("H40").Value = 0.482
$ws.Range("I40").Value = 0.5
$ws.Range("D41").Value = 0.057
$ws.Range("E41").Value = 0.691
$ws.Range("F41").Value = 0.069
$ws.Range("G41").Value = 0.106
$ws.Range("H41").Value = 0.049
$ws.Range("I41").Value = 0.094
$ws.Range("E42").Value = 0.913
$ws.Range("F42").Value = 0.001
$ws.Range("I42").Value = 0.002
$ws.Range("E43").Value = 0.989
